# Update excess mortality plots
# Applies the underlying data corrections for weeks 41-46/2022 (rows 131-151)
# and appends the newly-reported week-46 row (row 152) with its derived
# percentage-change formulas, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected raw counts for already-present weeks (rows 131-151) ---
$ws.Range("X131").Value = 708
$ws.Range("W133").Value = 473
$ws.Range("X135").Value = 641
$ws.Range("W142").Value = 458
$ws.Range("Y142").Value = 80
$ws.Range("X143").Value = 608
$ws.Range("R144").Value = 110
$ws.Range("U145").Value = 437
$ws.Range("W145").Value = 486
$ws.Range("Z145").Value = 433
$ws.Range("W146").Value = 502
$ws.Range("X146").Value = 664
$ws.Range("AA146").Value = 224
$ws.Range("V147").Value = 218
$ws.Range("V148").Value = 238
$ws.Range("W148").Value = 500
$ws.Range("X148").Value = 661
$ws.Range("Z148").Value = 535
$ws.Range("W149").Value = 515
$ws.Range("X149").Value = 641
$ws.Range("Y149").Value = 80
$ws.Range("Z149").Value = 499
$ws.Range("S150").Value = 228
$ws.Range("U150").Value = 389
$ws.Range("V150").Value = 234
$ws.Range("W150").Value = 478
$ws.Range("X150").Value = 656
$ws.Range("Z150").Value = 460
$ws.Range("P151").Value = 115
$ws.Range("Q151").Value = 137
$ws.Range("R151").Value = 124
$ws.Range("S151").Value = 214
$ws.Range("T151").Value = 62
$ws.Range("U151").Value = 414
$ws.Range("W151").Value = 438
$ws.Range("X151").Value = 607
$ws.Range("Y151").Value = 73
$ws.Range("Z151").Value = 474
$ws.Range("AA151").Value = 245

# --- New row 152: 2022 week 46 ---
$ws.Range("N152").Value = 2022
$ws.Range("O152").Value = 46
$ws.Range("P152").Value = 107
$ws.Range("Q152").Value = 150
$ws.Range("R152").Value = 94
$ws.Range("S152").Value = 228
$ws.Range("T152").Value = 57
$ws.Range("U152").Value = 426
$ws.Range("V152").Value = 240
$ws.Range("W152").Value = 498
$ws.Range("X152").Value = 651
$ws.Range("Y152").Value = 82
$ws.Range("Z152").Value = 512
$ws.Range("AA152").Value = 257
$ws.Range("AC152").Value = 2022
$ws.Range("AD152").Value = 46

# Percent-change formulas for the new row, matching the pattern used by
# the rows immediately above (e.g. row 151).
$ws.Range("AE152").Formula = "=ROUND((P152-B152)/B152*100,2)"
$ws.Range("AF152").Formula = "=ROUND((Q152-C152)/C152*100,2)"
$ws.Range("AG152").Formula = "=ROUND((R152-D152)/D152*100,2)"
$ws.Range("AH152").Formula = "=ROUND((S152-E152)/E152*100,2)"
$ws.Range("AI152").Formula = "=ROUND((T152-F152)/F152*100,2)"
$ws.Range("AJ152").Formula = "=ROUND((U152-G152)/G152*100,2)"
$ws.Range("AK152").Formula = "=ROUND((V152-H152)/H152*100,2)"
$ws.Range("AL152").Formula = "=ROUND((W152-I152)/I152*100,2)"
$ws.Range("AM152").Formula = "=ROUND((X152-J152)/J152*100,2)"
$ws.Range("AN152").Formula = "=ROUND((Y152-K152)/K152*100,2)"
$ws.Range("AO152").Formula = "=ROUND((Z152-L152)/L152*100,2)"
$ws.Range("AP152").Formula = "=ROUND((AA152-M152)/M152*100,2)"

# --- View-state refresh (zoom + selection), matching the updated sheet view ---
$ws.Range("AA99").Select()
$excel.ActiveWindow.Zoom = 70
